# UPDATE IN EXCEL READER, CONFIG FILE, CONFIG PROP
#
# Adds three new worksheets (Program, Batch, Sheet3) between the existing
# "Login" and "PythonArray" sheets, populates the new "Program" sheet with
# LMS program/config data, and restores the usual cursor/active-sheet state.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Insert "Program" right after "Login" and fill it with config data.
# ---------------------------------------------------------------------
$loginSheet = $wb.Worksheets.Item("Login")
$programSheet = $wb.Worksheets.Add($null, $loginSheet)
$programSheet.Name = "Program"

$programSheet.Range("A1").Value = "programDescription"
$programSheet.Range("A1").Font.Name = "Aptos"

$programSheet.Range("B1").Value = "programName"
$programSheet.Range("B1").Font.Name = "Calibri"

$programSheet.Range("C1").Value = "programStatus"
$programSheet.Range("C1").Font.Name = "Calibri"

$programSheet.Range("D1").Value = "programId"

$programSheet.Range("A2").Value = "Selenium Classes"
$programSheet.Range("A2").Font.Name = "Calibri"

$programSheet.Range("B2").Value = "ABCD1"
$programSheet.Range("B2").Font.Name = "Calibri"

$programSheet.Range("C2").Value = "Active"
$programSheet.Range("C2").Font.Name = "Calibri"

$programSheet.Range("D2").Value = 0

# ---------------------------------------------------------------------
# 2. Insert empty "Batch" sheet after "Program".
# ---------------------------------------------------------------------
$batchSheet = $wb.Worksheets.Add($null, $programSheet)
$batchSheet.Name = "Batch"

# ---------------------------------------------------------------------
# 3. Insert empty "Sheet3" after "Batch" (still before "PythonArray").
# ---------------------------------------------------------------------
$sheet3 = $wb.Worksheets.Add($null, $batchSheet)
$sheet3.Name = "Sheet3"

# ---------------------------------------------------------------------
# 4. Restore cursor positions / active sheet to match the saved view.
# ---------------------------------------------------------------------
$loginSheet.Activate() | Out-Null
$loginSheet.Range("D1").Select() | Out-Null

$programSheet.Activate() | Out-Null
$programSheet.Range("B2").Select() | Out-Null

$batchSheet.Activate() | Out-Null
$batchSheet.Range("E5").Select() | Out-Null

$programSheet.Activate() | Out-Null
